# Adds a new "2022-Q4" sheet (with fund-holdings detail) right after "总计"
# and before "2022-Q3", and inserts a matching summary row into "总计".

$wb = $excel.ActiveWorkbook

# Helper: write a piece of text into a cell while guaranteeing it is stored
# as TEXT (even when it looks numeric, e.g. "50.10") and that no stray
# number-format / quote-prefix style is left behind on the cell.
function Set-TextCell {
    param($cell, [string]$text)
    $cell.Value2 = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Duplicate an existing quarter sheet (2022-Q2 has 5 data rows, the
#    closest match to the 7 we need) right after "总计" so the new sheet
#    inherits identical sheet-level formatting/layout, then rename it.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsSource = $wb.Worksheets.Item("2022-Q2")
$wsSource.Copy($null, $wsTotal)

$newSheet = $wb.Worksheets.Item($wsTotal.Index + 1)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------
# 2) Overwrite the fund table with the 2022-Q4 data (7 funds).
# ---------------------------------------------------------------------
$rows = @(
    @("001481", "华宝油气（QDII）美元", "50.10", "94.65", "2.08", "1.0421", 10),
    @("007844", "华宝油气（QDII）人民币 C", "27.91", "94.65", "2.08", "0.5805", 10),
    @("162411", "华宝油气（QDII）人民币A", "22.19", "94.65", "2.08", "0.4616", 10),
    @("160416", "华安标普全球石油指数（QDII-LOF）A", "2.81", "93.63", "9.51", "0.2672", 2),
    @("014982", "华安标普全球石油指数（QDII-LOF）C", "0.36", "93.63", "9.51", "0.0342", 2),
    @("519981", "长信美国标准普尔100等权重指数增强（QDII）人民币", "0.44", "82.94", "0.87", "0.0038", 2),
    @("011706", "长信美国标准普尔100等权重指数增强（QDII）美元", "0.44", "82.94", "0.87", "0.0038", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2

    # Column A: the style-bearing numeric row index (same style as A2,
    # which already exists in the copied sheet).
    $newSheet.Cells.Item(2, 1).Copy()
    $newSheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $newSheet.Cells.Item($r, 1).Value2 = $i

    Set-TextCell $newSheet.Cells.Item($r, 2) $rows[$i][0]
    Set-TextCell $newSheet.Cells.Item($r, 3) $rows[$i][1]
    Set-TextCell $newSheet.Cells.Item($r, 4) $rows[$i][2]
    Set-TextCell $newSheet.Cells.Item($r, 5) $rows[$i][3]
    Set-TextCell $newSheet.Cells.Item($r, 6) $rows[$i][4]
    Set-TextCell $newSheet.Cells.Item($r, 7) $rows[$i][5]
    $newSheet.Cells.Item($r, 8).Value2 = $rows[$i][6]
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Insert a new "2022-Q4" row into "总计", shifting the existing
#    quarters down by one (a new "2020-Q4" row appears at the bottom).
# ---------------------------------------------------------------------
$totalRows = @(
    @("2022-Q4", 7, 2.39),
    @("2022-Q3", 2, 0.26),
    @("2022-Q2", 5, 2.6),
    @("2022-Q1", 2, 0.31),
    @("2021-Q4", 4, 1.88),
    @("2021-Q3", 1, 0.32),
    @("2021-Q2", 2, 0.45),
    @("2021-Q1", 5, 3.11),
    @("2020-Q4", 2, 0.93)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2

    $wsTotal.Cells.Item(2, 1).Copy()
    $wsTotal.Cells.Item($r, 1).PasteSpecial(-4122)
    $wsTotal.Cells.Item($r, 1).Value2 = $i

    $wsTotal.Cells.Item($r, 2).Value2 = $totalRows[$i][0]
    $wsTotal.Cells.Item($r, 3).Value2 = $totalRows[$i][1]
    $wsTotal.Cells.Item($r, 4).Value2 = $totalRows[$i][2]
}

$excel.CutCopyMode = 0

# Select 总计 as the active sheet isn't required, but put focus back on
# the last sheet (2020-Q4), matching the original tabSelected sheet.
$wb.Worksheets.Item("2020-Q4").Select()
